$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the hyperlinks that point at the (soon to be deleted) second
# resume row (D3: mailto link, E3: LinkedIn link). We look them up by
# their address instead of calling Range("D3").Hyperlinks.Delete() /
# Range("E3").Hyperlinks.Delete() because deleting through a Range's
# Hyperlinks collection removes *every* hyperlink on the sheet rather
# than just the one anchored at that range.
function Remove-HyperlinkAt($address) {
    foreach ($hl in $ws.Hyperlinks) {
        $hlAddress = $hl.Range.Address()
        if ($hlAddress -eq $address) {
            $hl.Delete()
            return
        }
    }
}

Remove-HyperlinkAt('$D$3')
Remove-HyperlinkAt('$E$3')

# Wipe out the second resume entry's data (row 3) while keeping the
# formatting (cell styles) that was already applied to it.
$ws.Range("A3:AM3").ClearContents()

# The row had a custom height sized for the wrapped text that used to
# live in it; let it size back down to the sheet's default height now
# that the row is empty.
$ws.Rows.Item(3).AutoFit()

# Update the view: scroll back to the left edge of the sheet and move
# the active selection onto the now-empty row.
$ws.Range("AN3").Select()
